$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999438891241,
    0.99649015012897113,
    0.98505608408413203,
    0.98247735996586716,
    0.97005262133322256,
    0.93986314136779048,
    0.93595448263506376,
    0.92979565801345365,
    0.92831233285649128,
    0.92777023833632311,
    0.9277525702070466,
    0.92797556845376439,
    0.93158975123293741,
    0.93425433970908633,
    0.93793043478215621,
    0.93542413826560611,
    0.93171639158989206,
    0.93060748218962352,
    0.99380161739892992,
    0.98668452166179088,
    0.98528602506200913,
    0.9840215163693018,
    0.97802131116924074,
    0.9650004606690592,
    0.95854346779685962,
    0.93620735629690044,
    0.93136043113986078,
    0.90988635647320804,
    0.89461449568573448,
    0.88804397906373245,
    0.88039006446011414,
    0.87871072584088439,
    0.87819071274502081
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
